# feat: add 2022-Q1 data
#
# Inserts a new worksheet "2022-Q1" (holding per-fund holding detail,
# identical column layout to "2021-Q4") positioned between the existing
# "2021-Q4" and "总计" sheets, and prepends a corresponding summary row
# to the "总计" sheet.

$wb = $excel.ActiveWorkbook
$wsQ4 = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Add($null, $wsQ4)
$wsQ1.Name = "2022-Q1"

# NOTE: fetch "总计" only *after* the insert above - Worksheets.Item()
# resolves by current tab position, so a reference captured before the
# Add() would silently end up pointing at the new sheet instead.
$wsTotal = $wb.Worksheets.Item("总计")

# Reuse the existing bold/centered/bordered header style from 2021-Q4!B1
# across the whole header row, and the bold/centered/bordered index-column
# style from 2021-Q4!A2 down column A.
$wsQ4.Range("B1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)
$wsQ4.Range("A2").Copy()
$wsQ1.Range("A2:A10").PasteSpecial(-4122)

$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# Helper: write a value into a cell as genuine text (preserves leading
# zeroes in fund codes, matches source export format for numeric-looking
# strings like "13.68").
function Set-TextCell($ws, $addr, $val) {
  $ws.Range($addr).NumberFormat = "@"
  $ws.Range($addr).Value = $val
}

$rows = @(
  @{ idx=0; code="290011"; name="泰信中小盘精选混合";                         scale="13.68"; pos="94.64"; pct="9.86"; mv="1.3488"; mvNum=$null; rank=4  },
  @{ idx=1; code="001970"; name="泰信鑫选灵活配置混合A";                       scale="3.03";  pos="93.96"; pct="9.89"; mv="0.2997"; mvNum=$null; rank=2  },
  @{ idx=2; code="002580"; name="泰信鑫选灵活配置混合C";                       scale="2.04";  pos="93.96"; pct="9.89"; mv="0.2018"; mvNum=$null; rank=2  },
  @{ idx=3; code="006279"; name="中金瑞祥灵活配置混合A";                       scale="2.10";  pos="59.54"; pct="3.06"; mv="0.0643"; mvNum=$null; rank=10 },
  @{ idx=4; code="011030"; name="达诚价值先锋灵活配置混合型证券投资基金A";     scale="0.47";  pos="80.83"; pct="6.11"; mv="0.0287"; mvNum=$null; rank=9  },
  @{ idx=5; code="011031"; name="达诚价值先锋灵活配置混合型证券投资基金C";     scale="0.25";  pos="80.83"; pct="6.11"; mv="0.0153"; mvNum=$null; rank=9  },
  @{ idx=6; code="002292"; name="诺安益鑫灵活配置混合";                       scale="0.30";  pos="50.08"; pct="4.10"; mv="0.0123"; mvNum=$null; rank=7  },
  @{ idx=7; code="002810"; name="金信转型创新成长灵活配置混合";               scale="0.18";  pos="81.12"; pct="3.04"; mv="0.0055"; mvNum=$null; rank=10 },
  @{ idx=8; code="006280"; name="中金瑞祥灵活配置混合C";                       scale="0.00";  pos="59.54"; pct="3.06"; mv=$null;    mvNum=0;    rank=10 }
)

$r = 2
foreach ($row in $rows) {
  $wsQ1.Range("A$r").Value = $row.idx
  Set-TextCell $wsQ1 "B$r" $row.code
  Set-TextCell $wsQ1 "C$r" $row.name
  Set-TextCell $wsQ1 "D$r" $row.scale
  Set-TextCell $wsQ1 "E$r" $row.pos
  Set-TextCell $wsQ1 "F$r" $row.pct
  if ($row.mv -ne $null) {
    Set-TextCell $wsQ1 "G$r" $row.mv
  } else {
    $wsQ1.Range("G$r").Value = $row.mvNum
  }
  $wsQ1.Range("H$r").Value = $row.rank
  $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" summary row to "总计" (existing "2021-Q4" row
#    shifts from row 2 down to row 3)
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()

# The freshly-inserted row inherits stray formatting from the row above;
# strip it back to the workbook default before writing the new values.
$wsTotal.Range("B2:D2").ClearFormats()

# A2 keeps the bold/centered/bordered index-column style used throughout
# the workbook (same style as 2021-Q4!A2 / 总计!A3).
$wsQ4.Range("A2").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
Set-TextCell $wsTotal "B2" "2022-Q1"
$wsTotal.Range("C2").Value = 9
$wsTotal.Range("D2").Value = 1.98

# Keep the shifted-down original row's index in sequence.
$wsTotal.Range("A3").Value = 1
